$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cVals = New-Object 'object[,]' 24,1
$cVals[0,0] = 0.04945618783204964
$cVals[1,0] = 0.0438850139643705
$cVals[2,0] = 0.04048237419669931
$cVals[3,0] = 0.0391002077597733
$cVals[4,0] = 0.03887096485854613
$cVals[5,0] = 0.04046371599599752
$cVals[6,0] = 0.04753142539337318
$cVals[7,0] = 0.06154038195178657
$cVals[8,0] = 0.07193269707904903
$cVals[9,0] = 0.07668415196170031
$cVals[10,0] = 0.07848699084442501
$cVals[11,0] = 0.07809855664281429
$cVals[12,0] = 0.07683240041694717
$cVals[13,0] = 0.07605731127773652
$cVals[14,0] = 0.07162267056303051
$cVals[15,0] = 0.06890838075879913
$cVals[16,0] = 0.06734943875028421
$cVals[17,0] = 0.06682198973675213
$cVals[18,0] = 0.06919708823096471
$cVals[19,0] = 0.07720420341543388
$cVals[20,0] = 0.08245816051922361
$cVals[21,0] = 0.07965207606893898
$cVals[22,0] = 0.06906655879708978
$cVals[23,0] = 0.05773369380345628
$ws.Range("C2:C25").Value = $cVals

$dVals = New-Object 'object[,]' 24,1
$dVals[0,0] = 0.01281438384477074
$dVals[1,0] = 0.01314579327325305
$dVals[2,0] = 0.01335878948337132
$dVals[3,0] = 0.01344795762859796
$dVals[4,0] = 0.01346290650629722
$dVals[5,0] = 0.01335998246591519
$dVals[6,0] = 0.01292666746974636
$dVals[7,0] = 0.01215333635368321
$dVals[8,0] = 0.01163306499416272
$dVals[9,0] = 0.0114070796191541
$dVals[10,0] = 0.01132306862161503
$dVals[11,0] = 0.01134109184752763
$dVals[12,0] = 0.01140013646197513
$dVals[13,0] = 0.01143650760984904
$dVals[14,0] = 0.01164805110277367
$dVals[15,0] = 0.01178058324328557
$dVals[16,0] = 0.0118578165761436
$dVals[17,0] = 0.01188413812379174
$dVals[18,0] = 0.01176637082687026
$dVals[19,0] = 0.01138275093348184
$dVals[20,0] = 0.01114116733255521
$dVals[21,0] = 0.01126925985276106
$dVals[22,0] = 0.01177279302621326
$dVals[23,0] = 0.01235423079047848
$ws.Range("D2:D25").Value = $dVals

$eVals = New-Object 'object[,]' 24,1
$eVals[0,0] = 0.4185757265276635
$eVals[1,0] = 0.3648985179177231
$eVals[2,0] = 0.3320673950023547
$eVals[3,0] = 0.3187176278844532
$eVals[4,0] = 0.3165025926020775
$eVals[5,0] = 0.3318872407732982
$eVals[6,0] = 0.4000398216122392
$eVals[7,0] = 0.534823712998687
$eVals[8,0] = 0.6347326044212735
$eVals[9,0] = 0.6804177274038636
$eVals[10,0] = 0.6977547522600531
$eVals[11,0] = 0.6940192153360556
$eVals[12,0] = 0.6818432940955148
$eVals[13,0] = 0.6743901104634489
$eVals[14,0] = 0.6317519894448367
$eVals[15,0] = 0.6056575445321357
$eVals[16,0] = 0.5906707402795917
$eVals[17,0] = 0.5856001555993373
$eVals[18,0] = 0.6084330403007101
$eVals[19,0] = 0.6854186264217788
$eVals[20,0] = 0.7359506573641852
$eVals[21,0] = 0.7089598170280311
$eVals[22,0] = 0.6071781914303642
$eVals[23,0] = 0.4982199127470608
$ws.Range("E2:E25").Value = $eVals

$fVals = New-Object 'object[,]' 24,1
$fVals[0,0] = 0.9170649819931782
$fVals[1,0] = 0.890610193046129
$fVals[2,0] = 0.8751832053469855
$fVals[3,0] = 0.869100001828869
$fVals[4,0] = 0.8681021224418828
$fVals[5,0] = 0.8751003439656273
$fVals[6,0] = 0.9077728000861214
$fVals[7,0] = 0.9784097686320052
$fVals[8,0] = 1.034440151703649
$fVals[9,0] = 1.060856226402365
$fVals[10,0] = 1.070994939921022
$fVals[11,0] = 1.068805324985831
$fVals[12,0] = 1.06168761612237
$fVals[13,0] = 1.057345526409577
$fVals[14,0] = 1.03273264612919
$fVals[15,0] = 1.017872525343023
$fVals[16,0] = 1.009412555020177
$fVals[17,0] = 1.006563062739588
$fVals[18,0] = 1.019445372588336
$fVals[19,0] = 1.063774564633022
$fVals[20,0] = 1.093537176254031
$fVals[21,0] = 1.07757922776122
$fVals[22,0] = 1.018734029014496
$fVals[23,0] = 0.9585842403963909
$ws.Range("F2:F25").Value = $fVals

$gVals = New-Object 'object[,]' 24,1
$gVals[0,0] = 0.7871259100182328
$gVals[1,0] = 0.75839825486824
$gVals[2,0] = 0.7415356845845196
$gVals[3,0] = 0.7348570102240899
$gVals[4,0] = 0.7337596101891535
$gVals[5,0] = 0.7414448353380862
$gVals[6,0] = 0.7770581160534391
$gVals[7,0] = 0.8531621004884755
$gVals[8,0] = 0.913050797251401
$gVals[9,0] = 0.9411930679733587
$gVals[10,0] = 0.951981717745042
$gVals[11,0] = 0.9496522867275985
$gVals[12,0] = 0.9420780007574194
$gVals[13,0] = 0.9374557705544362
$gVals[14,0] = 0.91122993435377
$gVals[15,0] = 0.8953732015755236
$gVals[16,0] = 0.8863372536937106
$gVals[17,0] = 0.8832922595154571
$gVals[18,0] = 0.8970524216655349
$gVals[19,0] = 0.9442991550065472
$gVals[20,0] = 0.9759469699515932
$gVals[21,0] = 0.9589846846224361
$gVals[22,0] = 0.8962929963429076
$gVals[23,0] = 0.8318871844046782
$ws.Range("G2:G25").Value = $gVals

$hVals = New-Object 'object[,]' 24,1
$hVals[0,0] = 0.7791840144811317
$hVals[1,0] = 0.7724216530243098
$hVals[2,0] = 0.7688684784993143
$hVals[3,0] = 0.7675702349455946
$hVals[4,0] = 0.7673636780524902
$hVals[5,0] = 0.7688503649962399
$hVals[6,0] = 0.7767274898261576
$hVals[7,0] = 0.7969699327830995
$hVals[8,0] = 0.8148284049363781
$hVals[9,0] = 0.8236151129985387
$hVals[10,0] = 0.8270388240760553
$hVals[11,0] = 0.82629716316805
$hVals[12,0] = 0.823894846102263
$hVals[13,0] = 0.8224359403600943
$hVals[14,0] = 0.8142676008507692
$hVals[15,0] = 0.8094271011666763
$hVals[16,0] = 0.8067052949300546
$hVals[17,0] = 0.8057944096028677
$hVals[18,0] = 0.8099359229520644
$hVals[19,0] = 0.8245978406818324
$hVals[20,0] = 0.8347425073698673
$hVals[21,0] = 0.8292763050485235
$hVals[22,0] = 0.8097056945102281
$hVals[23,0] = 0.7909739811219083
$ws.Range("H2:H25").Value = $hVals

$kVals = New-Object 'object[,]' 24,1
$kVals[0,0] = 1.5934994602502
$kVals[1,0] = 1.404877773744147
$kVals[2,0] = 1.289279243980332
$kVals[3,0] = 1.242224547857518
$kVals[4,0] = 1.234414283694662
$kVals[5,0] = 1.288644436887239
$kVals[6,0] = 1.528416812317516
$kVals[7,0] = 2.000423834841285
$kVals[8,0] = 2.348490100951778
$kVals[9,0] = 2.507154237842144
$kVals[10,0] = 2.567285917556262
$kVals[11,0] = 2.554333272132965
$kVals[12,0] = 2.512100313371491
$kVals[13,0] = 2.486237849956126
$kVals[14,0] = 2.338127775616897
$kVals[15,0] = 2.247352197521025
$kVals[16,0] = 2.19517114605236
$kVals[17,0] = 2.177508745179864
$kVals[18,0] = 2.25701222540539
$kVals[19,0] = 2.52450381068769
$kVals[20,0] = 2.699611691115308
$kVals[21,0] = 2.606126415355618
$kVals[22,0] = 2.252644905027239
$kVals[23,0] = 1.872521181384002
$ws.Range("K2:K25").Value = $kVals

$nVals = New-Object 'object[,]' 24,1
$nVals[0,0] = 0.9727383235607192
$nVals[1,0] = 0.9881805320442432
$nVals[2,0] = 0.9981615608558769
$nVals[3,0] = 1.002354299348678
$nVals[4,0] = 1.003058067759529
$nVals[5,0] = 0.9982175981796075
$nVals[6,0] = 0.9779590133703913
$nVals[7,0] = 0.9422037520638114
$nVals[8,0] = 0.9183674873758711
$nVals[9,0] = 0.9080548089233886
$nVals[10,0] = 0.9042262328997914
$nVals[11,0] = 0.9050473727400004
$nVals[12,0] = 0.9077382923993866
$nVals[13,0] = 0.9093965451007371
$nVals[14,0] = 0.9190521482542984
$nVals[15,0] = 0.9251116591115505
$nVals[16,0] = 0.9286468483848473
$nVals[17,0] = 0.9298523668424927
$nVals[18,0] = 0.9244614451567976
$nVals[19,0] = 0.9069458222598392
$nVals[20,0] = 0.8959450564223417
$nVals[21,0] = 0.9017753904242767
$nVals[22,0] = 0.9247552465121238
$nVals[23,0] = 0.9727383235607192
$ws.Range("N2:N25").Value = $nVals
